$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update stats for 2025-08 (row 21)
$ws.Range("B21").Value = 6226
$ws.Range("C21").Value = 983
$ws.Range("D21").Value = 5600378
$ws.Range("E21").Value = 899.5146161259236
$ws.Range("F21").Value = 8.071515361916326
$ws.Range("G21").Value = 3.691983122362874
$ws.Range("H21").Value = 27.83414274275206
